# Re-generated aggregated "charging station" rows: the B/C/D/E columns hold
# stringified Python lists (models / start SoC / end SoC / power) and F holds
# the numeric total power (sum of the power list) for that hour. Rows not
# listed below are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 2).Value = '[''Fiat 500 E'']'
$ws.Cells.Item(7, 3).Value = '[0.35]'
$ws.Cells.Item(7, 4).Value = '[0.9500000000000003]'
$ws.Cells.Item(7, 5).Value = '[14.400000000000007]'
$ws.Cells.Item(7, 6).Value = 14.40000000000001

$ws.Cells.Item(8, 2).Value = '[''Others'', ''VW ID.4'']'
$ws.Cells.Item(8, 3).Value = '[0.35, 0.35]'
$ws.Cells.Item(8, 4).Value = '[0.8000000000000002, 0.9000000000000002]'
$ws.Cells.Item(8, 5).Value = '[23.683750000000007, 42.35000000000002]'
$ws.Cells.Item(8, 6).Value = 66.03375000000003

$ws.Cells.Item(9, 2).Value = '[''Others'', ''Others'']'
$ws.Cells.Item(9, 3).Value = '[0.05, 0.1]'
$ws.Cells.Item(9, 4).Value = '[1.0, 0.9000000000000002]'
$ws.Cells.Item(9, 5).Value = '[49.999027777777776, 42.104444444444454]'
$ws.Cells.Item(9, 6).Value = 92.10347222222222

$ws.Cells.Item(11, 2).Value = '[]'
$ws.Cells.Item(11, 3).Value = '[]'
$ws.Cells.Item(11, 4).Value = '[]'
$ws.Cells.Item(11, 5).Value = '[]'
$ws.Cells.Item(11, 6).Value = 0

$ws.Cells.Item(12, 2).Value = '[]'
$ws.Cells.Item(12, 3).Value = '[]'
$ws.Cells.Item(12, 4).Value = '[]'
$ws.Cells.Item(12, 5).Value = '[]'
$ws.Cells.Item(12, 6).Value = 0

$ws.Cells.Item(13, 2).Value = '[''VW ID.3'', ''Others'', ''TESLA MODEL Y'', ''Smart FORTWO'']'
$ws.Cells.Item(13, 3).Value = '[0.4, 0.15, 0.2, 0.4]'
$ws.Cells.Item(13, 4).Value = '[0.6, 0.9500000000000003, 0.65, 0.9500000000000003]'
$ws.Cells.Item(13, 5).Value = '[11.599999999999998, 42.104444444444454, 33.75, 9.680000000000005]'
$ws.Cells.Item(13, 6).Value = 97.13444444444445

$ws.Cells.Item(14, 2).Value = '[''Fiat 500 E'']'
$ws.Cells.Item(14, 3).Value = '[0.05]'
$ws.Cells.Item(14, 4).Value = '[0.9000000000000002]'
$ws.Cells.Item(14, 5).Value = '[20.400000000000006]'
$ws.Cells.Item(14, 6).Value = 20.40000000000001

$ws.Cells.Item(15, 2).Value = '[''Tesla MODEL 3'']'
$ws.Cells.Item(15, 3).Value = '[0.25]'
$ws.Cells.Item(15, 4).Value = '[0.9500000000000003]'
$ws.Cells.Item(15, 5).Value = '[35.000000000000014]'
$ws.Cells.Item(15, 6).Value = 35.00000000000001

$ws.Cells.Item(16, 2).Value = '[''Others'', ''Others'']'
$ws.Cells.Item(16, 3).Value = '[0.35, 0.1]'
$ws.Cells.Item(16, 4).Value = '[0.9500000000000003, 0.8000000000000002]'
$ws.Cells.Item(16, 5).Value = '[31.578333333333347, 36.84138888888889]'
$ws.Cells.Item(16, 6).Value = 68.41972222222225

$ws.Cells.Item(17, 2).Value = '[''VW ID.5'']'
$ws.Cells.Item(17, 3).Value = '[0.2]'
$ws.Cells.Item(17, 4).Value = '[1.0]'
$ws.Cells.Item(17, 5).Value = '[61.6]'
$ws.Cells.Item(17, 6).Value = 61.6

$ws.Cells.Item(18, 2).Value = '[''Others'']'
$ws.Cells.Item(18, 3).Value = '[0.2]'
$ws.Cells.Item(18, 4).Value = '[0.7500000000000001]'
$ws.Cells.Item(18, 5).Value = '[28.946805555555557]'
$ws.Cells.Item(18, 6).Value = 28.94680555555556

$ws.Cells.Item(32, 2).Value = '[''Others'', ''Others'', ''MINI Cooper SE'']'
$ws.Cells.Item(32, 3).Value = '[0.3, 0.1, 0.2]'
$ws.Cells.Item(32, 4).Value = '[0.65, 0.9000000000000002, 0.9000000000000002]'
$ws.Cells.Item(32, 5).Value = '[18.420694444444447, 42.104444444444454, 20.230000000000004]'
$ws.Cells.Item(32, 6).Value = 80.75513888888891

$ws.Cells.Item(33, 2).Value = '[''Others'', ''Others'', ''MINI Cooper SE'']'
$ws.Cells.Item(33, 3).Value = '[0.35, 0.15, 0.2]'
$ws.Cells.Item(33, 4).Value = '[1.0, 0.9000000000000002, 0.9000000000000002]'
$ws.Cells.Item(33, 5).Value = '[34.20986111111111, 39.47291666666668, 20.230000000000004]'
$ws.Cells.Item(33, 6).Value = 93.91277777777779

$ws.Cells.Item(36, 2).Value = '[''Others'', ''Hyundai KONA 64 kWh'']'
$ws.Cells.Item(36, 3).Value = '[0.25, 0.45]'
$ws.Cells.Item(36, 4).Value = '[0.8000000000000002, 0.7500000000000001]'
$ws.Cells.Item(36, 5).Value = '[28.946805555555564, 19.200000000000006]'
$ws.Cells.Item(36, 6).Value = 48.14680555555557

$ws.Cells.Item(37, 2).Value = '[''Tesla MODEL 3'', ''Fiat 500 E'', ''Others'']'
$ws.Cells.Item(37, 3).Value = '[0.25, 0.15, 0.2]'
$ws.Cells.Item(37, 4).Value = '[0.7000000000000001, 1.0, 1.0]'
$ws.Cells.Item(37, 5).Value = '[22.500000000000004, 20.4, 42.10444444444445]'
$ws.Cells.Item(37, 6).Value = 85.00444444444446

$ws.Cells.Item(40, 2).Value = '[''MINI Cooper SE'', ''Tesla MODEL 3'', ''MINI Cooper SE'']'
$ws.Cells.Item(40, 3).Value = '[0.15, 0.1, 0.2]'
$ws.Cells.Item(40, 4).Value = '[0.8500000000000002, 0.9000000000000002, 0.9000000000000002]'
$ws.Cells.Item(40, 5).Value = '[20.230000000000004, 40.000000000000014, 20.230000000000004]'
$ws.Cells.Item(40, 6).Value = 80.46000000000002

$ws.Cells.Item(41, 2).Value = '[''MINI Cooper SE'', ''Others'']'
$ws.Cells.Item(41, 3).Value = '[0.1, 0.1]'
$ws.Cells.Item(41, 4).Value = '[0.8500000000000002, 0.8500000000000002]'
$ws.Cells.Item(41, 5).Value = '[21.675000000000004, 39.47291666666668]'
$ws.Cells.Item(41, 6).Value = 61.14791666666668
